$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths: column A a bit wider than the rest, columns B.. default-ish width
$ws.Columns("A:A").ColumnWidth = 9.6666666666667
$ws.Columns("B:AMK").ColumnWidth = 7.6666666666667

# Re-apply the Normal style to the cells that only ever carried the plain
# default formatting (header row + Volume/Change/Gain/Loss/Avg Gain/Avg
# Loss/RS/RSI columns). This mirrors the touch-up the sheet received where
# those cells picked up an explicit (but visually identical) style entry,
# while the Date column and the OHLC columns keep their own number formats.
$ws.Range("A1:N1").Style = "Normal"
$ws.Range("G2").Style = "Normal"
for ($r = 3; $r -le 15; $r++) {
    $ws.Range("G$r`:J$r").Style = "Normal"
}
$ws.Range("G16:N16").Style = "Normal"
$ws.Range("I17:J17").Style = "Normal"
